# Weekly update: insert the latest "Pepino dulce" price-report week at the
# top of the data block (row 177) and push the rest of the historical rows
# down by two. This mirrors how a new reporting week gets prepended to the
# daily/weekly consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 177:260 down to 179:262, inserting two blank rows.
$ws.Rows.Item(177).Resize(2).Insert()

# New row 177 - "Especial" no longer reported this week; first new record is "Primera".
$ws.Range("A177").Value = 12
$ws.Range("B177").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C177").Value = "Metropolitana"
$ws.Range("D177").Value = 45001
$ws.Range("E177").Value = 13
$ws.Range("F177").Value = 100112043
$ws.Range("G177").Value = "Pepino dulce"
$ws.Range("H177").Value = "Cultivar IV Región"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 150
$ws.Range("K177").Value = 15000
$ws.Range("L177").Value = 15000
$ws.Range("M177").Value = 15000
$ws.Range("N177").Value = "$/bandeja 18 kilos"
$ws.Range("O177").Value = "Provincia de Limarí"
$ws.Range("P177").Value = 833
$ws.Range("Q177").Value = 18
$ws.Range("R177").Value = "Hortaliza"

# New row 178 - "Segunda" quality for the same new reporting week.
$ws.Range("A178").Value = 12
$ws.Range("B178").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C178").Value = "Metropolitana"
$ws.Range("D178").Value = 45001
$ws.Range("E178").Value = 13
$ws.Range("F178").Value = 100112043
$ws.Range("G178").Value = "Pepino dulce"
$ws.Range("H178").Value = "Cultivar IV Región"
$ws.Range("I178").Value = "Segunda"
$ws.Range("J178").Value = 210
$ws.Range("K178").Value = 11000
$ws.Range("L178").Value = 11000
$ws.Range("M178").Value = 11000
$ws.Range("N178").Value = "$/bandeja 18 kilos"
$ws.Range("O178").Value = "Provincia de Limarí"
$ws.Range("P178").Value = 611
$ws.Range("Q178").Value = 18
$ws.Range("R178").Value = "Hortaliza"
